$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Capitalize connector words (de/del/la/las/los/el/y) in state/municipality names ---
$ws.Range("B5").Value = "Rincón De Romos"
$ws.Range("B26").Value = "Amatenango De La Frontera"
$ws.Range("B27").Value = "Amatenango Del Valle"
$ws.Range("B30").Value = "Bejucal De Ocampo"
$ws.Range("B32").Value = "Benemérito De Las Américas"
$ws.Range("B42").Value = "Comitán De Domínguez"
$ws.Range("B63").Value = "Marqués De Comillas"
$ws.Range("B64").Value = "Mazapa De Madero"
$ws.Range("B69").Value = "Ocozocoautla De Espinosa"
$ws.Range("B78").Value = "Salto De Agua"
$ws.Range("B79").Value = "San Cristóbal De Las Casas"
$ws.Range("B109").Value = "Guadalupe Y Calvo"
$ws.Range("B110").Value = "Hidalgo Del Parral"
$ws.Range("B128").Value = "San Juan De Sabinas"
$ws.Range("A137").Value = "Ciudad De México"
$ws.Range("A169").Value = "Estado De México"
$ws.Range("B169").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B171").Value = "Almoloya De Alquisiras"
$ws.Range("B172").Value = "Almoloya De Juárez"
$ws.Range("B177").Value = "Atizapán De Zaragoza"
$ws.Range("B182").Value = "Chapa De Mota"
$ws.Range("B185").Value = "Coacalco De Berriozábal"
$ws.Range("B189").Value = "Ecatepec De Morelos"
$ws.Range("B195").Value = "Ixtapan De La Sal"
$ws.Range("B196").Value = "Ixtapan Del Oro"
$ws.Range("B209").Value = "Naucalpan De Juárez"
$ws.Range("B217").Value = "San Felipe Del Progreso"
$ws.Range("B218").Value = "San Simón De Guerrero"
$ws.Range("B220").Value = "Soyaniquilpan De Juárez"
$ws.Range("B228").Value = "Tenango Del Valle"
$ws.Range("B234").Value = "Tlalnepantla De Baz"
$ws.Range("B239").Value = "Valle De Bravo"
$ws.Range("B240").Value = "Valle De Chalco Solidaridad"
$ws.Range("B241").Value = "Villa De Allende"
$ws.Range("B242").Value = "Villa Del Carbón"
$ws.Range("B254").Value = "Apaseo El Alto"
$ws.Range("B255").Value = "Apaseo El Grande"
$ws.Range("B262").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B274").Value = "San Diego De La Unión"
$ws.Range("B276").Value = "San Francisco Del Rincón"
$ws.Range("B278").Value = "San Luis De La Paz"
$ws.Range("B279").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B280").Value = "Silao De La Victoria"
$ws.Range("B283").Value = "Valle De Santiago"
$ws.Range("B289").Value = "Acapulco De Juárez"
$ws.Range("B292").Value = "Ajuchitlán Del Progreso"
$ws.Range("B293").Value = "Alcozauca De Guerrero"
$ws.Range("B297").Value = "Atenango Del Río"
$ws.Range("B298").Value = "Atlamajalcingo Del Monte"
$ws.Range("B300").Value = "Atoyac De Álvarez"
$ws.Range("B301").Value = "Ayutla De Los Libres"
$ws.Range("B304").Value = "Buenavista De Cuéllar"
$ws.Range("B305").Value = "Chilapa De Álvarez"
$ws.Range("B306").Value = "Chilpancingo De Los Bravo"
$ws.Range("B307").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B312").Value = "Coyuca De Benítez"
$ws.Range("B313").Value = "Coyuca De Catalán"
$ws.Range("B317").Value = "Cutzamala De Pinzón"
$ws.Range("B323").Value = "Huitzuco De Los Figueroa"
$ws.Range("B324").Value = "Iguala De La Independencia"
$ws.Range("B326").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B327").Value = "Zihuatanejo De Azueta"
$ws.Range("B329").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B332").Value = "Mártir De Cuilapan"
$ws.Range("B343").Value = "Taxco De Alarcón"
$ws.Range("B345").Value = "Técpan De Galeana"
$ws.Range("B347").Value = "Tepecoacuilco De Trujano"
$ws.Range("B348").Value = "Tixtla De Guerrero"
$ws.Range("B351").Value = "Tlapa De Comonfort"
$ws.Range("B365").Value = "Atotonilco De Tula"
$ws.Range("B366").Value = "Atotonilco El Grande"
$ws.Range("B372").Value = "Cuautepec De Hinojosa"
$ws.Range("B375").Value = "Huasca De Ocampo"
$ws.Range("B379").Value = "Huejutla De Reyes"
$ws.Range("B382").Value = "Jacala De Ledezma"
$ws.Range("B388").Value = "Mineral De La Reforma"
$ws.Range("B389").Value = "Mineral Del Monte"
$ws.Range("B390").Value = "Mixquiahuala De Juárez"
$ws.Range("B391").Value = "Molango De Escamilla"
$ws.Range("B393").Value = "Nopala De Villagrán"
$ws.Range("B394").Value = "Pachuca De Soto"
$ws.Range("B397").Value = "Progreso De Obregón"
$ws.Range("B402").Value = "Santiago De Anaya"
$ws.Range("B403").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B407").Value = "Tenango De Doria"
$ws.Range("B409").Value = "Tepehuacán De Guerrero"
$ws.Range("B410").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B412").Value = "Tezontepec De Aldama"
$ws.Range("B418").Value = "Tula De Allende"
$ws.Range("B419").Value = "Tulancingo De Bravo"
$ws.Range("B423").Value = "Zacualtipán De Ángeles"
$ws.Range("B427").Value = "Ahualulco De Mercado"
$ws.Range("B429").Value = "Atotonilco El Alto"
$ws.Range("B435").Value = "Encarnación De Díaz"
$ws.Range("B440").Value = "Ixtlahuacán Del Río"
$ws.Range("B442").Value = "Jilotlán De Los Dolores"
$ws.Range("B444").Value = "Lagos De Moreno"
$ws.Range("B446").Value = "Ojuelos De Jalisco"
$ws.Range("B451").Value = "San Juan De Los Lagos"
$ws.Range("B453").Value = "San Miguel El Alto"
$ws.Range("B455").Value = "Talpa De Allende"
$ws.Range("B458").Value = "Teocuitatlán De Corona"
$ws.Range("B459").Value = "Tepatitlán De Morelos"
$ws.Range("B461").Value = "Tizapán El Alto"
$ws.Range("B462").Value = "Tlajomulco De Zúñiga"
$ws.Range("B471").Value = "Zapotlán El Grande"
$ws.Range("B490").Value = "Cojumatlán De Régules"
$ws.Range("B532").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B550").Value = "Coatlán Del Río"
$ws.Range("B561").Value = "Puente De Ixtla"
$ws.Range("B565").Value = "Tetela Del Volcán"
$ws.Range("B566").Value = "Tlaltizapán De Zapata"
$ws.Range("B574").Value = "Bahía De Banderas"
$ws.Range("B577").Value = "Ixtlán Del Río"
$ws.Range("B581").Value = "Santa María Del Oro"
$ws.Range("B595").Value = "Lampazos De Naranjo"
$ws.Range("B598").Value = "San Nicolás De Los Garza"
$ws.Range("B603").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B604").Value = "Ayoquezco De Aldama"
$ws.Range("B606").Value = "Chalcatongo De Hidalgo"
$ws.Range("B607").Value = "Ciénega De Zimatlán"
$ws.Range("B610").Value = "Constancia Del Rosario"
$ws.Range("B612").Value = "Fresnillo De Trujano"
$ws.Range("B613").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B614").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B615").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B616").Value = "Huautla De Jiménez"
$ws.Range("B617").Value = "Ixtlán De Juárez"
$ws.Range("B618").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B627").Value = "Mártires De Tacubaya"
$ws.Range("B629").Value = "Mazatlán Villa De Flores"
$ws.Range("B631").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B632").Value = "Mixistlán De La Reforma"
$ws.Range("B635").Value = "Oaxaca De Juárez"
$ws.Range("B636").Value = "Ocotlán De Morelos"
$ws.Range("B637").Value = "Pinotepa De Don Luis"
$ws.Range("B639").Value = "Putla Villa De Guerrero"
$ws.Range("B650").Value = "San Antonino El Alto"
$ws.Range("B659").Value = "San Dionisio Del Mar"
$ws.Range("B662").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B676").Value = "San José Del Progreso"
$ws.Range("B681").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B689").Value = "San Juan De Los Cués"
$ws.Range("B690").Value = "San Juan Del Estado"
$ws.Range("B719").Value = "San Miguel Del Puerto"
$ws.Range("B720").Value = "San Miguel El Grande"
$ws.Range("B727").Value = "San Pablo Villa De Mitla"
$ws.Range("B744").Value = "Santa Ana Del Valle"
$ws.Range("B750").Value = "Santa Cruz De Bravo"
$ws.Range("B762").Value = "Santa María Del Tule"
$ws.Range("B766").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B796").Value = "Santo Domingo De Morelos"
$ws.Range("B807").Value = "Tataltepec De Valdés"
$ws.Range("B808").Value = "Teotitlán De Flores Magón"
$ws.Range("B810").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B811").Value = "Tlacolula De Matamoros"
$ws.Range("B812").Value = "Tlalixtac De Cabrera"
$ws.Range("B813").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B814").Value = "Villa Sola De Vega"
$ws.Range("B815").Value = "Villa Talea De Castro"
$ws.Range("B816").Value = "Zimatlán De Álvarez"
$ws.Range("B829").Value = "Ayotoxco De Guerrero"
$ws.Range("B852").Value = "Ixcamilpa De Guerrero"
$ws.Range("B854").Value = "Izúcar De Matamoros"
$ws.Range("B863").Value = "Palmar De Bravo"
$ws.Range("B870").Value = "San Diego La Mesa Tochimiltzingo"
$ws.Range("B878").Value = "San Salvador El Seco"
$ws.Range("B879").Value = "San Salvador El Verde"
$ws.Range("B886").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B889").Value = "Tepexi De Rodríguez"
$ws.Range("B890").Value = "Tetela De Ocampo"
$ws.Range("B894").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B906").Value = "Xayacatlán De Bravo"
$ws.Range("B909").Value = "Xochitlán De Vicente Suárez"
$ws.Range("B919").Value = "Amealco De Bonfil"
$ws.Range("B921").Value = "Cadereyta De Montes"
$ws.Range("B926").Value = "Jalpan De Serra"
$ws.Range("B927").Value = "Landa De Matamoros"
$ws.Range("B929").Value = "Pinal De Amoles"
$ws.Range("B932").Value = "San Juan Del Río"
$ws.Range("B943").Value = "Armadillo De Los Infante"
$ws.Range("B944").Value = "Axtla De Terrazas"
$ws.Range("B949").Value = "Ciudad Del Maíz"
$ws.Range("B958").Value = "Mexquitic De Carmona"
$ws.Range("B962").Value = "San Ciro De Acosta"
$ws.Range("B968").Value = "Santa María Del Río"
$ws.Range("B969").Value = "Soledad De Graciano Sánchez"
$ws.Range("B978").Value = "Villa De Arriaga"
$ws.Range("B979").Value = "Villa De La Paz"
$ws.Range("B980").Value = "Villa De Reyes"
$ws.Range("B1016").Value = "Jalpa De Méndez"
$ws.Range("B1037").Value = "Soto La Marina"
$ws.Range("B1048").Value = "Contla De Juan Cuamatzi"
$ws.Range("B1051").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B1052").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B1055").Value = "Papalotla De Xicohténcatl"
$ws.Range("B1057").Value = "San Pablo Del Monte"
$ws.Range("B1058").Value = "Sanctórum De Lázaro Cárdenas"
$ws.Range("B1075").Value = "Amatlán De Los Reyes"
$ws.Range("B1082").Value = "Boca Del Río"
$ws.Range("B1084").Value = "Camarón De Tejeda"
$ws.Range("B1088").Value = "Castillo De Teayo"
$ws.Range("B1105").Value = "Cosamaloapan De Carpio"
$ws.Range("B1106").Value = "Cosautlán De Carvajal"
$ws.Range("B1123").Value = "Hueyapan De Ocampo"
$ws.Range("B1124").Value = "Ignacio De La Llave"
$ws.Range("B1127").Value = "Ixhuacán De Los Reyes"
$ws.Range("B1128").Value = "Ixhuatlán De Madero"
$ws.Range("B1129").Value = "Ixhuatlán Del Café"
$ws.Range("B1140").Value = "Juchique De Ferrer"
$ws.Range("B1143").Value = "Lerdo De Tejada"
$ws.Range("B1146").Value = "Martínez De La Torre"
$ws.Range("B1151").Value = "Mixtla De Altamirano"
$ws.Range("B1153").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B1160").Value = "Ozuluama De Mascareñas"
$ws.Range("B1163").Value = "Paso De Ovejas"
$ws.Range("B1164").Value = "Paso Del Macho"
$ws.Range("B1168").Value = "Poza Rica De Hidalgo"
$ws.Range("B1176").Value = "Sayula De Alemán"
$ws.Range("B1179").Value = "Soledad De Doblado"
$ws.Range("B1182").Value = "Tatahuicapan De Juárez"
$ws.Range("B1199").Value = "Tlacotepec De Mejía"
$ws.Range("B1210").Value = "Vega De Alatorre"
$ws.Range("B1220").Value = "Zozocolco De Hidalgo"
$ws.Range("B1233").Value = "Concepción Del Oro"
$ws.Range("B1239").Value = "Jiménez Del Teul"
$ws.Range("B1242").Value = "Nochistlán De Mejía"
$ws.Range("B1243").Value = "Noria De Ángeles"

# --- Minor floating point recalculation artifacts on two Total rows ---
$ws.Range("D102").Value = 0.09185480301018148
$ws.Range("D1256").Value = 0.009185480301018147

# --- Remove trailing footnote rows (1258-1263) and shrink dimension to A1:D1257 ---
$ws.Range("A1258:A1263").EntireRow.Delete()
